# Auto-generated edit script applying the "Updated cryptos list" diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking Price (D) cells that must stay as literal text
# (otherwise Excel auto-converts strings like "0.163" into the number 0.163,
# losing the fixed-width "51.768.34"-style grouping used by this sheet).
$textCells = @("D5","D6","D7","D9","D10","D12","D13","D14","D16","D19","D20","D21","D23","D24","D25","D26","D27","D28","D29","D30","D32","D34","D35","D36","D38","D39","D40","D43","D44","D46","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Cell value updates
$ws.Range("D2").Value = '51.711.55'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '2.824.18'
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '350.61'
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("D6").Value = '112.82'
$ws.Range("E6").Value = '  +4.54%  '
$ws.Range("D7").Value = '0.557'
$ws.Range("E7").Value = '  +1.49%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.618'
$ws.Range("E9").Value = '  +5.42%  '
$ws.Range("D10").Value = '40.10'
$ws.Range("E10").Value = '  +1.36%  '
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("D12").Value = '0.0847'
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").Value = '19.94'
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("D14").Value = '7.76'
$ws.Range("E14").Value = '  +2.87%  '
$ws.Range("D15").Value = '3.273.01'
$ws.Range("E15").Value = '  +1.94%  '
$ws.Range("D16").Value = '0.972'
$ws.Range("E16").Value = '  +5.73%  '
$ws.Range("D17").Value = '2.821.36'
$ws.Range("E17").Value = '  +1.74%  '
$ws.Range("D18").Value = '51.747.65'
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("D19").Value = '3.47'
$ws.Range("E19").Value = '  +12.10%  '
$ws.Range("D20").Value = '7.58'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '13.38'
$ws.Range("E21").Value = '  +2.00%  '
$ws.Range("D22").Value = '0.0₃0971'
$ws.Range("E22").Value = '  +0.91%  '
$ws.Range("D23").Value = '70.47'
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("D24").Value = '268.13'
$ws.Range("E24").Value = '  +1.10%  '
$ws.Range("D25").Value = '2.75'
$ws.Range("E25").Value = '  +1.88%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '26.19'
$ws.Range("E27").Value = '  +0.78%  '
$ws.Range("D28").Value = '0.163'
$ws.Range("E28").Value = '  +0.33%  '
$ws.Range("D29").Value = '10.58'
$ws.Range("E29").Value = '  +3.62%  '
$ws.Range("D30").Value = '38.83'
$ws.Range("E30").Value = '  +6.33%  '
$ws.Range("E31").Value = '  +2.84%  '
$ws.Range("D32").Value = '6.27'
$ws.Range("E32").Value = '  +1.95%  '
$ws.Range("E33").Value = '  +1.76%  '
$ws.Range("D34").Value = '0.0896'
$ws.Range("E34").Value = '  +8.34%  '
$ws.Range("D35").Value = '0.0451'
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").Value = '5.64'
$ws.Range("E36").Value = '  +1.95%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").Value = '19.01'
$ws.Range("E38").Value = '  +3.67%  '
$ws.Range("D39").Value = '3.21'
$ws.Range("E39").Value = '  +2.24%  '
$ws.Range("D40").Value = '2.01'
$ws.Range("E40").Value = '  +2.51%  '
$ws.Range("E41").Value = '  +1.79%  '
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("D43").Value = '122.72'
$ws.Range("E43").Value = '  +2.35%  '
$ws.Range("D44").Value = '22.16'
$ws.Range("E44").Value = '  +0.75%  '
$ws.Range("E45").Value = '  +1.06%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '2.51'
$ws.Range("E46").Value = '  +8.60%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.179.10'
$ws.Range("E47").Value = '  +3.67%  '
$ws.Range("D48").Value = '3.49'
$ws.Range("E48").Value = '  +7.62%  '
$ws.Range("D49").Value = '0.247'
$ws.Range("E49").Value = '  +23.34%  '
$ws.Range("D50").Value = '0.951'
$ws.Range("E50").Value = '  +5.71%  '
$ws.Range("E51").Value = '  +2.18%  '
